# Updating with latest games: add games 21-24 (Final Fantasy x2, Doctor Who, Takir Dragonstorm)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (centered style) from the last existing data rows down onto
# the new block so the new cells pick up the same style index as the rest of
# the table, without introducing any new style entries.
$ws.Range("A60:G75").Copy()
$ws.Range("A76:G91").PasteSpecial(-4122)

$ws.Range("A76").Value = 21
$ws.Range("B76").Value = 'Alex'
$ws.Range("C76").Value = 'Final Fantasy'
$ws.Range("D76").Value = 'Limit Breaker'
$ws.Range("E76").Value = 2
$ws.Range("F76").Value = 'Houston'
$ws.Range("G76").Value = 'Red'
$ws.Range("A77").Value = 21
$ws.Range("B77").Value = 'Chris'
$ws.Range("C77").Value = 'Final Fantasy'
$ws.Range("D77").Value = 'Scions & Spell Craft'
$ws.Range("E77").Value = 3
$ws.Range("F77").Value = 'Houston'
$ws.Range("G77").Value = 'White'
$ws.Range("A78").Value = 21
$ws.Range("B78").Value = 'Kevin'
$ws.Range("C78").Value = 'Final Fantasy'
$ws.Range("D78").Value = 'Counter Blitz'
$ws.Range("E78").Value = 1
$ws.Range("F78").Value = 'Houston'
$ws.Range("G78").Value = 'Green'
$ws.Range("A79").Value = 21
$ws.Range("B79").Value = 'Sandro'
$ws.Range("C79").Value = 'Final Fantasy'
$ws.Range("D79").Value = 'Revival Trance'
$ws.Range("E79").Value = 4
$ws.Range("F79").Value = 'Houston'
$ws.Range("G79").Value = 'Red'
$ws.Range("A80").Value = 22
$ws.Range("B80").Value = 'Alex'
$ws.Range("C80").Value = 'Final Fantasy'
$ws.Range("D80").Value = 'Scions & Spell Craft'
$ws.Range("E80").Value = 2
$ws.Range("F80").Value = 'Houston'
$ws.Range("G80").Value = 'White'
$ws.Range("A81").Value = 22
$ws.Range("B81").Value = 'Chris'
$ws.Range("C81").Value = 'Final Fantasy'
$ws.Range("D81").Value = 'Revival Trance'
$ws.Range("E81").Value = 1
$ws.Range("F81").Value = 'Houston'
$ws.Range("G81").Value = 'Red'
$ws.Range("A82").Value = 22
$ws.Range("B82").Value = 'Kevin'
$ws.Range("C82").Value = 'Final Fantasy'
$ws.Range("D82").Value = 'Limit Breaker'
$ws.Range("E82").Value = 3
$ws.Range("F82").Value = 'Houston'
$ws.Range("G82").Value = 'Red'
$ws.Range("A83").Value = 22
$ws.Range("B83").Value = 'Sandro'
$ws.Range("C83").Value = 'Final Fantasy'
$ws.Range("D83").Value = 'Counter Blitz'
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = 'Houston'
$ws.Range("G83").Value = 'Green'
$ws.Range("A84").Value = 23
$ws.Range("B84").Value = 'Alex'
$ws.Range("C84").Value = 'Doctor Who'
$ws.Range("D84").Value = 'Timey Wimey'
$ws.Range("E84").Value = 2
$ws.Range("F84").Value = 'Houston'
$ws.Range("G84").Value = 'Blue'
$ws.Range("A85").Value = 23
$ws.Range("B85").Value = 'Chris'
$ws.Range("C85").Value = 'Doctor Who'
$ws.Range("D85").Value = 'Masters of Evil'
$ws.Range("E85").Value = 3
$ws.Range("F85").Value = 'Houston'
$ws.Range("G85").Value = 'Black'
$ws.Range("A86").Value = 23
$ws.Range("B86").Value = 'Kevin'
$ws.Range("C86").Value = 'Doctor Who'
$ws.Range("D86").Value = 'Blast From the Past'
$ws.Range("E86").Value = 1
$ws.Range("F86").Value = 'Houston'
$ws.Range("G86").Value = 'Green'
$ws.Range("A87").Value = 23
$ws.Range("B87").Value = 'Sandro'
$ws.Range("C87").Value = 'Doctor Who'
$ws.Range("D87").Value = 'Paradox Power'
$ws.Range("E87").Value = 4
$ws.Range("F87").Value = 'Houston'
$ws.Range("G87").Value = 'Green'
$ws.Range("A88").Value = 24
$ws.Range("B88").Value = 'Alex'
$ws.Range("C88").Value = 'Takir Dragonstorm'
$ws.Range("D88").Value = 'Temur Roar'
$ws.Range("E88").Value = 1
$ws.Range("F88").Value = 'Houston'
$ws.Range("G88").Value = 'White'
$ws.Range("A89").Value = 24
$ws.Range("B89").Value = 'Chris'
$ws.Range("C89").Value = 'Takir Dragonstorm'
$ws.Range("D89").Value = 'Mardu Surge'
$ws.Range("E89").Value = 4
$ws.Range("F89").Value = 'Houston'
$ws.Range("G89").Value = 'Green'
$ws.Range("A90").Value = 24
$ws.Range("B90").Value = 'Kevin'
$ws.Range("C90").Value = 'Takir Dragonstorm'
$ws.Range("D90").Value = 'Saultai Anisen'
$ws.Range("E90").Value = 2
$ws.Range("F90").Value = 'Houston'
$ws.Range("G90").Value = 'Black'
$ws.Range("A91").Value = 24
$ws.Range("B91").Value = 'Sandro'
$ws.Range("C91").Value = 'Takir Dragonstorm'
$ws.Range("D91").Value = 'Jeski Striker'
$ws.Range("E91").Value = 3
$ws.Range("F91").Value = 'Houston'
$ws.Range("G91").Value = 'Red'

# Widen column D to fit the new, longer deck names. (COM ColumnWidth adds
# the standard ~0.83 char padding before it is stored as the OOXML <col>
# width, so 20.16 here round-trips to the target stored width of 21.)
$ws.Columns.Item(4).ColumnWidth = 20.16

# Re-anchor the selection to the top-left cell (closest achievable to "no
# explicit selection").
$ws.Range("A1").Select()

# Grow the AutoFilter to cover the newly added rows.
$ws.AutoFilterMode = $false
$ws.Range("A1:G91").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$91"
    }
}
